$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 21:46"

$data = New-Object 'object[,]' 184,8
$data[0,0] = "China"
$data[0,1] = 80967
$data[0,2] = 39
$data[0,3] = 71150
$data[0,4] = 6569
$data[0,5] = 2136
$data[0,6] = 3
$data[0,7] = 3248
$data[1,0] = "Italia"
$data[1,1] = 47021
$data[1,2] = 5986
$data[1,3] = 5129
$data[1,4] = 37860
$data[1,5] = 2655
$data[1,6] = 627
$data[1,7] = 4032
$data[2,0] = "España"
$data[2,1] = 20412
$data[2,2] = 2335
$data[2,3] = 1588
$data[2,4] = 17774
$data[2,5] = 939
$data[2,6] = 219
$data[2,7] = 1050
$data[3,0] = "Alemania"
$data[3,1] = 19848
$data[3,2] = 4528
$data[3,3] = 180
$data[3,4] = 19601
$data[3,5] = 2
$data[3,6] = 23
$data[3,7] = 67
$data[4,0] = "Iran"
$data[4,1] = 19644
$data[4,2] = 1237
$data[4,3] = 6745
$data[4,4] = 11466
$data[4,5] = 0
$data[4,6] = 149
$data[4,7] = 1433
$data[5,0] = "Estados Unidos"
$data[5,1] = 18121
$data[5,2] = 4332
$data[5,3] = 125
$data[5,4] = 17763
$data[5,5] = 64
$data[5,6] = 26
$data[5,7] = 233
$data[6,0] = "Francia"
$data[6,1] = 12612
$data[6,2] = 1617
$data[6,3] = 1295
$data[6,4] = 10867
$data[6,5] = 1122
$data[6,6] = 78
$data[6,7] = 450
$data[7,0] = "Corea del Sur"
$data[7,1] = 8652
$data[7,2] = 87
$data[7,3] = 2233
$data[7,4] = 6325
$data[7,5] = 59
$data[7,6] = 3
$data[7,7] = 94
$data[8,0] = "Suiza"
$data[8,1] = 5407
$data[8,2] = 1185
$data[8,3] = 15
$data[8,4] = 5336
$data[8,5] = 0
$data[8,6] = 13
$data[8,7] = 56
$data[9,0] = "Reino Unido"
$data[9,1] = 3983
$data[9,2] = 714
$data[9,3] = 65
$data[9,4] = 3741
$data[9,5] = 20
$data[9,6] = 33
$data[9,7] = 177
$data[10,0] = "Paises Bajos"
$data[10,1] = 2994
$data[10,2] = 534
$data[10,3] = 2
$data[10,4] = 2886
$data[10,5] = 210
$data[10,6] = 30
$data[10,7] = 106
$data[11,0] = "Austria"
$data[11,1] = 2491
$data[11,2] = 312
$data[11,3] = 9
$data[11,4] = 2476
$data[11,5] = 14
$data[11,6] = 0
$data[11,7] = 6
$data[12,0] = "Belgica"
$data[12,1] = 2257
$data[12,2] = 462
$data[12,3] = 204
$data[12,4] = 2016
$data[12,5] = 164
$data[12,6] = 16
$data[12,7] = 37
$data[13,0] = "Noruega"
$data[13,1] = 1926
$data[13,2] = 136
$data[13,3] = 1
$data[13,4] = 1918
$data[13,5] = 27
$data[13,6] = 0
$data[13,7] = 7
$data[14,0] = "Suecia"
$data[14,1] = 1639
$data[14,2] = 200
$data[14,3] = 16
$data[14,4] = 1607
$data[14,5] = 21
$data[14,6] = 5
$data[14,7] = 16
$data[15,0] = "Dinamarca"
$data[15,1] = 1255
$data[15,2] = 104
$data[15,3] = 1
$data[15,4] = 1245
$data[15,5] = 37
$data[15,6] = 3
$data[15,7] = 9
$data[16,0] = "Malasia"
$data[16,1] = 1030
$data[16,2] = 130
$data[16,3] = 87
$data[16,4] = 940
$data[16,5] = 26
$data[16,6] = 1
$data[16,7] = 3
$data[17,0] = "Portugal"
$data[17,1] = 1020
$data[17,2] = 234
$data[17,3] = 5
$data[17,4] = 1009
$data[17,5] = 26
$data[17,6] = 2
$data[17,7] = 6
$data[18,0] = "Japon"
$data[18,1] = 963
$data[18,2] = 20
$data[18,3] = 215
$data[18,4] = 715
$data[18,5] = 50
$data[18,6] = 0
$data[18,7] = 33
$data[19,0] = "Canada"
$data[19,1] = 943
$data[19,2] = 70
$data[19,3] = 11
$data[19,4] = 920
$data[19,5] = 1
$data[19,6] = 0
$data[19,7] = 12
$data[20,0] = "Australia"
$data[20,1] = 876
$data[20,2] = 120
$data[20,3] = 46
$data[20,4] = 823
$data[20,5] = 2
$data[20,6] = 0
$data[20,7] = 7
$data[21,0] = "Chequia"
$data[21,1] = 833
$data[21,2] = 139
$data[21,3] = 4
$data[21,4] = 829
$data[21,5] = 6
$data[21,6] = 0
$data[21,7] = 0
$data[22,0] = "Brasil"
$data[22,1] = 819
$data[22,2] = 179
$data[22,3] = 2
$data[22,4] = 806
$data[22,5] = 18
$data[22,6] = 4
$data[22,7] = 11
$data[23,0] = "Crucero"
$data[23,1] = 712
$data[23,2] = 0
$data[23,3] = 527
$data[23,4] = 177
$data[23,5] = 14
$data[23,6] = 1
$data[23,7] = 8
$data[24,0] = "Israel"
$data[24,1] = 705
$data[24,2] = 28
$data[24,3] = 15
$data[24,4] = 689
$data[24,5] = 10
$data[24,6] = 1
$data[24,7] = 1
$data[25,0] = "Irlanda"
$data[25,1] = 683
$data[25,2] = 126
$data[25,3] = 5
$data[25,4] = 675
$data[25,5] = 6
$data[25,6] = 0
$data[25,7] = 3
$data[26,0] = "Turquia"
$data[26,1] = 670
$data[26,2] = 311
$data[26,3] = 0
$data[26,4] = 666
$data[26,5] = 0
$data[26,6] = 0
$data[26,7] = 4
$data[27,0] = "Pakistan"
$data[27,1] = 500
$data[27,2] = 46
$data[27,3] = 13
$data[27,4] = 484
$data[27,5] = 0
$data[27,6] = 1
$data[27,7] = 3
$data[28,0] = "Grecia"
$data[28,1] = 495
$data[28,2] = 31
$data[28,3] = 19
$data[28,4] = 466
$data[28,5] = 20
$data[28,6] = 4
$data[28,7] = 10
$data[29,0] = "Luxemburgo"
$data[29,1] = 484
$data[29,2] = 149
$data[29,3] = 6
$data[29,4] = 473
$data[29,5] = 1
$data[29,6] = 1
$data[29,7] = 5
$data[30,0] = "Catar"
$data[30,1] = 470
$data[30,2] = 10
$data[30,3] = 10
$data[30,4] = 460
$data[30,5] = 6
$data[30,6] = 0
$data[30,7] = 0
$data[31,0] = "Finlandia"
$data[31,1] = 450
$data[31,2] = 50
$data[31,3] = 10
$data[31,4] = 440
$data[31,5] = 2
$data[31,6] = 0
$data[31,7] = 0
$data[32,0] = "Chile"
$data[32,1] = 434
$data[32,2] = 92
$data[32,3] = 6
$data[32,4] = 428
$data[32,5] = 7
$data[32,6] = 0
$data[32,7] = 0
$data[33,0] = "Polonia"
$data[33,1] = 425
$data[33,2] = 70
$data[33,3] = 13
$data[33,4] = 407
$data[33,5] = 3
$data[33,6] = 0
$data[33,7] = 5
$data[34,0] = "Islandia"
$data[34,1] = 409
$data[34,2] = 79
$data[34,3] = 5
$data[34,4] = 404
$data[34,5] = 1
$data[34,6] = 0
$data[34,7] = 0
$data[35,0] = "Singapur"
$data[35,1] = 385
$data[35,2] = 40
$data[35,3] = 131
$data[35,4] = 254
$data[35,5] = 14
$data[35,6] = 0
$data[35,7] = 0
$data[36,0] = "Indonesia"
$data[36,1] = 369
$data[36,2] = 60
$data[36,3] = 17
$data[36,4] = 320
$data[36,5] = 0
$data[36,6] = 7
$data[36,7] = 32
$data[37,0] = "Ecuador"
$data[37,1] = 367
$data[37,2] = 107
$data[37,3] = 3
$data[37,4] = 359
$data[37,5] = 2
$data[37,6] = 2
$data[37,7] = 5
$data[38,0] = "Arabia Saudita"
$data[38,1] = 344
$data[38,2] = 70
$data[38,3] = 8
$data[38,4] = 336
$data[38,5] = 0
$data[38,6] = 0
$data[38,7] = 0
$data[39,0] = "Eslovenia"
$data[39,1] = 341
$data[39,2] = 22
$data[39,3] = 0
$data[39,4] = 340
$data[39,5] = 6
$data[39,6] = 0
$data[39,7] = 1
$data[40,0] = "Tailandia"
$data[40,1] = 322
$data[40,2] = 50
$data[40,3] = 42
$data[40,4] = 279
$data[40,5] = 1
$data[40,6] = 0
$data[40,7] = 1
$data[41,0] = "Rumania"
$data[41,1] = 308
$data[41,2] = 31
$data[41,3] = 31
$data[41,4] = 277
$data[41,5] = 11
$data[41,6] = 0
$data[41,7] = 0
$data[42,0] = "Barein"
$data[42,1] = 297
$data[42,2] = 18
$data[42,3] = 125
$data[42,4] = 171
$data[42,5] = 4
$data[42,6] = 0
$data[42,7] = 1
$data[43,0] = "Egipto"
$data[43,1] = 285
$data[43,2] = 29
$data[43,3] = 42
$data[43,4] = 235
$data[43,5] = 0
$data[43,6] = 1
$data[43,7] = 8
$data[44,0] = "Estonia"
$data[44,1] = 283
$data[44,2] = 16
$data[44,3] = 1
$data[44,4] = 282
$data[44,5] = 1
$data[44,6] = 0
$data[44,7] = 0
$data[45,0] = "Peru"
$data[45,1] = 263
$data[45,2] = 29
$data[45,3] = 1
$data[45,4] = 258
$data[45,5] = 5
$data[45,6] = 3
$data[45,7] = 4
$data[46,0] = "Hong Kong"
$data[46,1] = 256
$data[46,2] = 48
$data[46,3] = 98
$data[46,4] = 154
$data[46,5] = 4
$data[46,6] = 0
$data[46,7] = 4
$data[47,0] = "Rusia"
$data[47,1] = 253
$data[47,2] = 54
$data[47,3] = 12
$data[47,4] = 240
$data[47,5] = 0
$data[47,6] = 0
$data[47,7] = 1
$data[48,0] = "India"
$data[48,1] = 249
$data[48,2] = 55
$data[48,3] = 23
$data[48,4] = 221
$data[48,5] = 0
$data[48,6] = 1
$data[48,7] = 5
$data[49,0] = "Filipinas"
$data[49,1] = 230
$data[49,2] = 13
$data[49,3] = 8
$data[49,4] = 204
$data[49,5] = 1
$data[49,6] = 1
$data[49,7] = 18
$data[50,0] = "Irak"
$data[50,1] = 208
$data[50,2] = 16
$data[50,3] = 49
$data[50,4] = 142
$data[50,5] = 0
$data[50,6] = 4
$data[50,7] = 17
$data[51,0] = "Sudafrica"
$data[51,1] = 202
$data[51,2] = 52
$data[51,3] = 0
$data[51,4] = 202
$data[51,5] = 0
$data[51,6] = 0
$data[51,7] = 0
$data[52,0] = "Libano"
$data[52,1] = 177
$data[52,2] = 20
$data[52,3] = 4
$data[52,4] = 169
$data[52,5] = 3
$data[52,6] = 0
$data[52,7] = 4
$data[53,0] = "Mexico"
$data[53,1] = 164
$data[53,2] = 46
$data[53,3] = 4
$data[53,4] = 159
$data[53,5] = 1
$data[53,6] = 0
$data[53,7] = 1
$data[54,0] = "Kuwait"
$data[54,1] = 159
$data[54,2] = 11
$data[54,3] = 22
$data[54,4] = 137
$data[54,5] = 5
$data[54,6] = 0
$data[54,7] = 0
$data[55,0] = "Colombia"
$data[55,1] = 145
$data[55,2] = 37
$data[55,3] = 1
$data[55,4] = 144
$data[55,5] = 0
$data[55,6] = 0
$data[55,7] = 0
$data[56,0] = "San Marino"
$data[56,1] = 144
$data[56,2] = 0
$data[56,3] = 4
$data[56,4] = 126
$data[56,5] = 12
$data[56,6] = 0
$data[56,7] = 14
$data[57,0] = "Emiratos Arabes Unidos"
$data[57,1] = 140
$data[57,2] = 0
$data[57,3] = 31
$data[57,4] = 109
$data[57,5] = 2
$data[57,6] = 0
$data[57,7] = 0
$data[58,0] = "Eslovaquia"
$data[58,1] = 137
$data[58,2] = 13
$data[58,3] = 0
$data[58,4] = 137
$data[58,5] = 2
$data[58,6] = 0
$data[58,7] = 0
$data[59,0] = "Panama"
$data[59,1] = 137
$data[59,2] = 0
$data[59,3] = 1
$data[59,4] = 135
$data[59,5] = 7
$data[59,6] = 0
$data[59,7] = 1
$data[60,0] = "Armenia"
$data[60,1] = 136
$data[60,2] = 14
$data[60,3] = 1
$data[60,4] = 135
$data[60,5] = 2
$data[60,6] = 0
$data[60,7] = 0
$data[61,0] = "Serbia"
$data[61,1] = 135
$data[61,2] = 32
$data[61,3] = 2
$data[61,4] = 132
$data[61,5] = 4
$data[61,6] = 1
$data[61,7] = 1
$data[62,0] = "Taiwan"
$data[62,1] = 135
$data[62,2] = 27
$data[62,3] = 28
$data[62,4] = 105
$data[62,5] = 0
$data[62,6] = 1
$data[62,7] = 2
$data[63,0] = "Bulgaria"
$data[63,1] = 129
$data[63,2] = 22
$data[63,3] = 1
$data[63,4] = 125
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 3
$data[64,0] = "Argentina"
$data[64,1] = 128
$data[64,2] = 0
$data[64,3] = 3
$data[64,4] = 122
$data[64,5] = 0
$data[64,6] = 0
$data[64,7] = 3
$data[65,0] = "Croacia"
$data[65,1] = 128
$data[65,2] = 18
$data[65,3] = 5
$data[65,4] = 122
$data[65,5] = 0
$data[65,6] = 0
$data[65,7] = 1
$data[66,0] = "Costa Rica"
$data[66,1] = 113
$data[66,2] = 26
$data[66,3] = 2
$data[66,4] = 109
$data[66,5] = 2
$data[66,6] = 1
$data[66,7] = 2
$data[67,0] = "Letonia"
$data[67,1] = 111
$data[67,2] = 25
$data[67,3] = 1
$data[67,4] = 110
$data[67,5] = 0
$data[67,6] = 0
$data[67,7] = 0
$data[68,0] = "Uruguay"
$data[68,1] = 94
$data[68,2] = 15
$data[68,3] = 0
$data[68,4] = 94
$data[68,5] = 0
$data[68,6] = 0
$data[68,7] = 0
$data[69,0] = "Vietnam"
$data[69,1] = 91
$data[69,2] = 6
$data[69,3] = 17
$data[69,4] = 74
$data[69,5] = 0
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = "Argelia"
$data[70,1] = 90
$data[70,2] = 0
$data[70,3] = 32
$data[70,4] = 47
$data[70,5] = 0
$data[70,6] = 2
$data[70,7] = 11
$data[71,0] = "Bosnia y Herzegovina"
$data[71,1] = 89
$data[71,2] = 25
$data[71,3] = 2
$data[71,4] = 87
$data[71,5] = 1
$data[71,6] = 0
$data[71,7] = 0
$data[72,0] = "Jordania"
$data[72,1] = 85
$data[72,2] = 16
$data[72,3] = 1
$data[72,4] = 84
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 0
$data[73,0] = "Hungria"
$data[73,1] = 85
$data[73,2] = 12
$data[73,3] = 7
$data[73,4] = 74
$data[73,5] = 6
$data[73,6] = 3
$data[73,7] = 4
$data[74,0] = "Islas Feroe"
$data[74,1] = 80
$data[74,2] = 8
$data[74,3] = 3
$data[74,4] = 77
$data[74,5] = 0
$data[74,6] = 0
$data[74,7] = 0
$data[75,0] = "Brunei"
$data[75,1] = 78
$data[75,2] = 5
$data[75,3] = 1
$data[75,4] = 77
$data[75,5] = 2
$data[75,6] = 0
$data[75,7] = 0
$data[76,0] = "Marruecos"
$data[76,1] = 77
$data[76,2] = 14
$data[76,3] = 2
$data[76,4] = 72
$data[76,5] = 1
$data[76,6] = 1
$data[76,7] = 3
$data[77,0] = "Republica de Chipre"
$data[77,1] = 75
$data[77,2] = 8
$data[77,3] = 0
$data[77,4] = 75
$data[77,5] = 1
$data[77,6] = 0
$data[77,7] = 0
$data[78,0] = "Principado de Andorra"
$data[78,1] = 75
$data[78,2] = 1
$data[78,3] = 1
$data[78,4] = 74
$data[78,5] = 2
$data[78,6] = 0
$data[78,7] = 0
$data[79,0] = "Sri Lanka"
$data[79,1] = 73
$data[79,2] = 13
$data[79,3] = 3
$data[79,4] = 70
$data[79,5] = 0
$data[79,6] = 0
$data[79,7] = 0
$data[80,0] = "Republica Dominicana"
$data[80,1] = 72
$data[80,2] = 38
$data[80,3] = 0
$data[80,4] = 70
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 2
$data[81,0] = "Republica de Macedonia"
$data[81,1] = 70
$data[81,2] = 20
$data[81,3] = 1
$data[81,4] = 69
$data[81,5] = 1
$data[81,6] = 0
$data[81,7] = 0
$data[82,0] = "Albania"
$data[82,1] = 70
$data[82,2] = 6
$data[82,3] = 0
$data[82,4] = 68
$data[82,5] = 2
$data[82,6] = 0
$data[82,7] = 2
$data[83,0] = "Bielorrusia"
$data[83,1] = 69
$data[83,2] = 18
$data[83,3] = 15
$data[83,4] = 54
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 0
$data[84,0] = "Moldavia"
$data[84,1] = 66
$data[84,2] = 17
$data[84,3] = 1
$data[84,4] = 64
$data[84,5] = 3
$data[84,6] = 0
$data[84,7] = 1
$data[85,0] = "Malta"
$data[85,1] = 64
$data[85,2] = 11
$data[85,3] = 2
$data[85,4] = 62
$data[85,5] = 1
$data[85,6] = 0
$data[85,7] = 0
$data[86,0] = "Lituania"
$data[86,1] = 63
$data[86,2] = 15
$data[86,3] = 1
$data[86,4] = 62
$data[86,5] = 1
$data[86,6] = 0
$data[86,7] = 0
$data[87,0] = "Tunez"
$data[87,1] = 54
$data[87,2] = 15
$data[87,3] = 1
$data[87,4] = 52
$data[87,5] = 7
$data[87,6] = 0
$data[87,7] = 1
$data[88,0] = "Kazajistan"
$data[88,1] = 52
$data[88,2] = 8
$data[88,3] = 0
$data[88,4] = 52
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 0
$data[89,0] = "Camboya"
$data[89,1] = 51
$data[89,2] = 14
$data[89,3] = 1
$data[89,4] = 50
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 0
$data[90,0] = "Oman"
$data[90,1] = 48
$data[90,2] = 0
$data[90,3] = 13
$data[90,4] = 35
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 0
$data[91,0] = "Estado de Palestina"
$data[91,1] = 48
$data[91,2] = 1
$data[91,3] = 17
$data[91,4] = 31
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 0
$data[92,0] = "Guadalupe"
$data[92,1] = 45
$data[92,2] = 12
$data[92,3] = 0
$data[92,4] = 45
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 0
$data[93,0] = "Georgia"
$data[93,1] = 44
$data[93,2] = 4
$data[93,3] = 1
$data[93,4] = 43
$data[93,5] = 1
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = "Azerbaiyan"
$data[94,1] = 44
$data[94,2] = 0
$data[94,3] = 7
$data[94,4] = 36
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 1
$data[95,0] = "Venezuela"
$data[95,1] = 42
$data[95,2] = 0
$data[95,3] = 0
$data[95,4] = 42
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 0
$data[96,0] = "Burkina Faso"
$data[96,1] = 40
$data[96,2] = 7
$data[96,3] = 4
$data[96,4] = 35
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 1
$data[97,0] = "Nueva Zelanda"
$data[97,1] = 39
$data[97,2] = 11
$data[97,3] = 0
$data[97,4] = 39
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 0
$data[98,0] = "Ucrania"
$data[98,1] = 39
$data[98,2] = 13
$data[98,3] = 1
$data[98,4] = 35
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 3
$data[99,0] = "Senegal"
$data[99,1] = 38
$data[99,2] = 2
$data[99,3] = 2
$data[99,4] = 36
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 0
$data[100,0] = "Uzbekistan"
$data[100,1] = 33
$data[100,2] = 10
$data[100,3] = 0
$data[100,4] = 33
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 0
$data[101,0] = "Martinica"
$data[101,1] = 32
$data[101,2] = 9
$data[101,3] = 0
$data[101,4] = 31
$data[101,5] = 7
$data[101,6] = 0
$data[101,7] = 1
$data[102,0] = "Reunion"
$data[102,1] = 28
$data[102,2] = 0
$data[102,3] = 0
$data[102,4] = 28
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 0
$data[103,0] = "Liechtenstein"
$data[103,1] = 28
$data[103,2] = 0
$data[103,3] = 0
$data[103,4] = 28
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 0
$data[104,0] = "Honduras"
$data[104,1] = 24
$data[104,2] = 12
$data[104,3] = 0
$data[104,4] = 24
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 0
$data[105,0] = "Afganistan"
$data[105,1] = 24
$data[105,2] = 2
$data[105,3] = 1
$data[105,4] = 23
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 0
$data[106,0] = "Camerun"
$data[106,1] = 20
$data[106,2] = 7
$data[106,3] = 2
$data[106,4] = 18
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 0
$data[107,0] = "Banglades"
$data[107,1] = 20
$data[107,2] = 2
$data[107,3] = 3
$data[107,4] = 16
$data[107,5] = 1
$data[107,6] = 0
$data[107,7] = 1
$data[108,0] = "Consejo Danes para los Refugiados"
$data[108,1] = 18
$data[108,2] = 4
$data[108,3] = 0
$data[108,4] = 18
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 0
$data[109,0] = "Ruanda"
$data[109,1] = 17
$data[109,2] = 6
$data[109,3] = 0
$data[109,4] = 17
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 0
$data[110,0] = "Macao"
$data[110,1] = 17
$data[110,2] = 0
$data[110,3] = 10
$data[110,4] = 7
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 0
$data[111,0] = "Bolivia"
$data[111,1] = 16
$data[111,2] = 1
$data[111,3] = 0
$data[111,4] = 16
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 0
$data[112,0] = "Ghana"
$data[112,1] = 16
$data[112,2] = 5
$data[112,3] = 0
$data[112,4] = 16
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 0
$data[113,0] = "Cuba"
$data[113,1] = 16
$data[113,2] = 5
$data[113,3] = 0
$data[113,4] = 15
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 1
$data[114,0] = "Jamaica"
$data[114,1] = 16
$data[114,2] = 1
$data[114,3] = 2
$data[114,4] = 13
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 1
$data[115,0] = "Guayana Francesa"
$data[115,1] = 15
$data[115,2] = 0
$data[115,3] = 0
$data[115,4] = 15
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 0
$data[116,0] = "Guyana"
$data[116,1] = 15
$data[116,2] = 10
$data[116,3] = 0
$data[116,4] = 14
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 1
$data[117,0] = "Guam"
$data[117,1] = 14
$data[117,2] = 2
$data[117,3] = 0
$data[117,4] = 14
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 0
$data[118,0] = "Montenegro"
$data[118,1] = 14
$data[118,2] = 1
$data[118,3] = 0
$data[118,4] = 14
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 0
$data[119,0] = "Puerto Rico"
$data[119,1] = 14
$data[119,2] = 8
$data[119,3] = 0
$data[119,4] = 14
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = "Maldivas"
$data[120,1] = 13
$data[120,2] = 0
$data[120,3] = 0
$data[120,4] = 13
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 0
$data[121,0] = "Paraguay"
$data[121,1] = 13
$data[121,2] = 0
$data[121,3] = 0
$data[121,4] = 13
$data[121,5] = 1
$data[121,6] = 0
$data[121,7] = 0
$data[122,0] = "Mauricio"
$data[122,1] = 12
$data[122,2] = 5
$data[122,3] = 0
$data[122,4] = 12
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = "Guatemala"
$data[123,1] = 12
$data[123,2] = 3
$data[123,3] = 0
$data[123,4] = 11
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 1
$data[124,0] = "Nigeria"
$data[124,1] = 12
$data[124,2] = 0
$data[124,3] = 1
$data[124,4] = 11
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = "Polinesia Francesa"
$data[125,1] = 11
$data[125,2] = 5
$data[125,3] = 0
$data[125,4] = 11
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = "Monaco"
$data[126,1] = 11
$data[126,2] = 1
$data[126,3] = 0
$data[126,4] = 11
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 0
$data[127,0] = "Gibraltar"
$data[127,1] = 10
$data[127,2] = 0
$data[127,3] = 2
$data[127,4] = 8
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = "Etiopia"
$data[128,1] = 9
$data[128,2] = 2
$data[128,3] = 0
$data[128,4] = 9
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = "Trinidad yTobago"
$data[129,1] = 9
$data[129,2] = 0
$data[129,3] = 0
$data[129,4] = 9
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = "Togo"
$data[130,1] = 9
$data[130,2] = 8
$data[130,3] = 0
$data[130,4] = 9
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = "Costa de Marfil"
$data[131,1] = 9
$data[131,2] = 0
$data[131,3] = 1
$data[131,4] = 8
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 0
$data[132,0] = "Kenia"
$data[132,1] = 7
$data[132,2] = 0
$data[132,3] = 0
$data[132,4] = 7
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = "Seychelles"
$data[133,1] = 7
$data[133,2] = 1
$data[133,3] = 0
$data[133,4] = 7
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = "Guinea Ecuatorial"
$data[134,1] = 6
$data[134,2] = 0
$data[134,3] = 0
$data[134,4] = 6
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = "Tanzania"
$data[135,1] = 6
$data[135,2] = 0
$data[135,3] = 0
$data[135,4] = 6
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = "Mongolia"
$data[136,1] = 6
$data[136,2] = 0
$data[136,3] = 0
$data[136,4] = 6
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = "Mayotte"
$data[137,1] = 6
$data[137,2] = 2
$data[137,3] = 0
$data[137,4] = 6
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 0
$data[138,0] = "Kirguistan"
$data[138,1] = 6
$data[138,2] = 3
$data[138,3] = 0
$data[138,4] = 6
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = "Barbados"
$data[139,1] = 5
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 5
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = "Aruba"
$data[140,1] = 5
$data[140,2] = 0
$data[140,3] = 1
$data[140,4] = 4
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = "Surinam"
$data[141,1] = 4
$data[141,2] = 3
$data[141,3] = 0
$data[141,4] = 4
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = "San Martin (Parte Francesa)"
$data[142,1] = 4
$data[142,2] = 1
$data[142,3] = 0
$data[142,4] = 4
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = "Bahamas"
$data[143,1] = 4
$data[143,2] = 1
$data[143,3] = 0
$data[143,4] = 4
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = "Congo"
$data[144,1] = 3
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 3
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = "Republica de Africa Central"
$data[145,1] = 3
$data[145,2] = 2
$data[145,3] = 0
$data[145,4] = 3
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = "Madagascar"
$data[146,1] = 3
$data[146,2] = 3
$data[146,3] = 0
$data[146,4] = 3
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = "Namibia"
$data[147,1] = 3
$data[147,2] = 0
$data[147,3] = 0
$data[147,4] = 3
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 0
$data[148,0] = "San Bartolome"
$data[148,1] = 3
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 3
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = "Islas Virgenes de los Estados Unidos"
$data[149,1] = 3
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 3
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = "Islas Caimanes"
$data[150,1] = 3
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 2
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 1
$data[151,0] = "Curazao"
$data[151,1] = 3
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 2
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 1
$data[152,0] = "Gabon"
$data[152,1] = 3
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 2
$data[152,5] = 0
$data[152,6] = 1
$data[152,7] = 1
$data[153,0] = "Nueva Caledonia"
$data[153,1] = 2
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 2
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = "Mauritania"
$data[154,1] = 2
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 2
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = "Guinea"
$data[155,1] = 2
$data[155,2] = 1
$data[155,3] = 0
$data[155,4] = 2
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = "Bermudas"
$data[156,1] = 2
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 2
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = "Santa Lucia"
$data[157,1] = 2
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 2
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = "Benin"
$data[158,1] = 2
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 2
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = "Groenlandia"
$data[159,1] = 2
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 2
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = "Isla de Man"
$data[160,1] = 2
$data[160,2] = 1
$data[160,3] = 0
$data[160,4] = 2
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = "Haiti"
$data[161,1] = 2
$data[161,2] = 2
$data[161,3] = 0
$data[161,4] = 2
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = "Zambia"
$data[162,1] = 2
$data[162,2] = 0
$data[162,3] = 0
$data[162,4] = 2
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = "Butan"
$data[163,1] = 2
$data[163,2] = 1
$data[163,3] = 0
$data[163,4] = 2
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = "Liberia"
$data[164,1] = 2
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 2
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = "Sudan"
$data[165,1] = 2
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 1
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 1
$data[166,0] = "El Salvador"
$data[166,1] = 1
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 1
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = "Fiyi"
$data[167,1] = 1
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 1
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = "Nicaragua"
$data[168,1] = 1
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 1
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = "Republica del Chad"
$data[169,1] = 1
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 1
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = "Niger"
$data[170,1] = 1
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 1
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = "Cabo Verde"
$data[171,1] = 1
$data[171,2] = 1
$data[171,3] = 0
$data[171,4] = 1
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = "Papua Nueva Guinea"
$data[172,1] = 1
$data[172,2] = 1
$data[172,3] = 0
$data[172,4] = 1
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = "Angola"
$data[173,1] = 1
$data[173,2] = 1
$data[173,3] = 0
$data[173,4] = 1
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = "Republica de Yibuti"
$data[174,1] = 1
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 1
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = "Suazilandia"
$data[175,1] = 1
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 1
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = "San Martin (Parte Holandesa)"
$data[176,1] = 1
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 1
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = "Santa Sede"
$data[177,1] = 1
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 1
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "San Vicente y las Granadinas"
$data[178,1] = 1
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 1
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = "Montserrat"
$data[179,1] = 1
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 1
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = "Somalia"
$data[180,1] = 1
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 1
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = "Antigua y Barbuda"
$data[181,1] = 1
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 1
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = "Gambia"
$data[182,1] = 1
$data[182,2] = 0
$data[182,3] = 0
$data[182,4] = 1
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = "Nepal"
$data[183,1] = 1
$data[183,2] = 0
$data[183,3] = 1
$data[183,4] = 0
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0

$ws.Range("A4:H187").Value = $data
